$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the boilerplate columns from the current row 319 (market, region,
# category id/name, variety, unit, kg/unidades flag, classification) before
# the insert shifts it down to row 320 - they are identical across this
# "Camote" block and carry over unchanged into the freshly inserted row.
$mercadoId  = $ws.Range("A319").Value2
$mercado    = $ws.Range("B319").Value2
$region     = $ws.Range("C319").Value2
$codreg     = $ws.Range("E319").Value2
$categoriaId= $ws.Range("F319").Value2
$categoria  = $ws.Range("G319").Value2
$variedad   = $ws.Range("H319").Value2
$unidad     = $ws.Range("N319").Value2
$kgUnidades = $ws.Range("Q319").Value2
$clasif     = $ws.Range("R319").Value2
$dateFmt    = $ws.Range("D319").NumberFormat

# Insert a new blank row at 319, pushing the existing rows 319-348 down to
# 320-349.
$ws.Rows.Item(319).Insert()

# Fill the newly inserted row 319 with its data.
$ws.Range("A319").Value = $mercadoId
$ws.Range("B319").Value = $mercado
$ws.Range("C319").Value = $region
$ws.Range("D319").Value = 45131
$ws.Range("D319").NumberFormat = $dateFmt
$ws.Range("E319").Value = $codreg
$ws.Range("F319").Value = $categoriaId
$ws.Range("G319").Value = $categoria
$ws.Range("H319").Value = $variedad
$ws.Range("I319").Value = "1a (guarda)"
$ws.Range("J319").Value = 280
$ws.Range("K319").Value = 550
$ws.Range("L319").Value = 550
$ws.Range("M319").Value = 550
$ws.Range("N319").Value = $unidad
$ws.Range("O319").Value = "Región del Maule"
$ws.Range("P319").Value = 550
$ws.Range("Q319").Value = $kgUnidades
$ws.Range("R319").Value = $clasif
